# PlayerPerformance_5393.xlsx update
# - Adds a new "Player Info" sheet (before "ODI Batting") with player bio data.
# - Adds a new "ODI Batting Extra" sheet (after "ODI Batting") with extra
#   per-match batting stats.
# - Renames "ODI Batting"'s MATCH_CARD_LINK column to MATCH_CODE and replaces
#   the full howstat.com scorecard URLs with the bare numeric match code.

$wb = $excel.ActiveWorkbook
$odiBatting = $wb.Worksheets.Item("ODI Batting")

# ---------------------------------------------------------------------------
# 1. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code
# ---------------------------------------------------------------------------
$matchCodes = @("4519", "4520", "4522", "4533", "4535", "4536", "4577", "4580", "4583", "4586", "4590", "4592", "4606", "4611", "4616", "4621", "4623", "4624", "4636", "4639", "4642", "4727", "4731")

$odiBatting.Cells.Item(1, 4).Value = "MATCH_CODE"
for ($i = 0; $i -lt $matchCodes.Count; $i++) {
    $row = $i + 2
    $odiBatting.Cells.Item($row, 4).NumberFormat = "@"
    $odiBatting.Cells.Item($row, 4).Value = $matchCodes[$i]
}

# ---------------------------------------------------------------------------
# 2. New "Player Info" sheet, placed before "ODI Batting"
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

$playerHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $playerHeaders.Count; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $playerHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$playerInfo.Cells.Item(2, 1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "5393"
$playerInfo.Cells.Item(2, 2).Value = "Shamarh Shaqad Joshua Brooks"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Leg Break"

# ---------------------------------------------------------------------------
# 3. New "ODI Batting Extra" sheet, placed after "ODI Batting"
# ---------------------------------------------------------------------------
# Re-fetch a fresh reference: the sheet collection shifted after the
# "Player Info" insert above and the old $odiBatting handle no longer
# anchors correctly.
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$extra = $wb.Worksheets.Add($null, $odiBatting)
$extra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Count; $c++) {
    $cell = $extra.Cells.Item(1, $c)
    $cell.Value = $extraHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$extraRows = @(
    @("4533", $null, $null, $null, $null, "NO"),
    @("4535", 4, "2", "2", "22.80%", "NO"),
    @("4536", 4, "0", "0", $null, "NO"),
    @("4577", 2, "3", "2", "24.10%", "NO"),
    @("4580", 2, "0", "0", "2.76%", "NO"),
    @("4583", 3, "3", "4", "32.79%", "NO"),
    @("4586", $null, $null, $null, $null, "NO"),
    @("4590", $null, $null, $null, $null, "NO"),
    @("4592", 3, "1", "1", "8.33%", "NO"),
    @("4606", $null, $null, $null, $null, "NO"),
    @("4611", 3, "1", "0", "4.63%", "NO"),
    @("4616", 3, "1", "0", "2.25%", "NO"),
    @("4621", 3, "4", "1", "15.08%", "NO"),
    @("4623", $null, $null, $null, $null, "NO"),
    @("4624", 3, "0", "0", $null, "NO"),
    @("4636", 3, "9", "1", "40.93%", "YES"),
    @("4639", 3, "0", "0", $null, "NO"),
    @("4642", $null, $null, $null, $null, "NO"),
    @("4727", $null, $null, $null, $null, "NO"),
    @("4731", 3, "2", "0", "6.92%", "NO")
)

for ($i = 0; $i -lt $extraRows.Count; $i++) {
    $row = $i + 2
    $rowData = $extraRows[$i]

    $extra.Cells.Item($row, 1).NumberFormat = "@"
    $extra.Cells.Item($row, 1).Value = $rowData[0]

    if ($null -ne $rowData[1]) {
        $extra.Cells.Item($row, 2).Value = $rowData[1]
    }
    if ($null -ne $rowData[2]) {
        $extra.Cells.Item($row, 3).NumberFormat = "@"
        $extra.Cells.Item($row, 3).Value = $rowData[2]
    }
    if ($null -ne $rowData[3]) {
        $extra.Cells.Item($row, 4).NumberFormat = "@"
        $extra.Cells.Item($row, 4).Value = $rowData[3]
    }
    if ($null -ne $rowData[4]) {
        $extra.Cells.Item($row, 5).NumberFormat = "@"
        $extra.Cells.Item($row, 5).Value = $rowData[4]
    }
    $extra.Cells.Item($row, 6).Value = $rowData[5]
}

# Keep the first sheet ("Player Info") active/selected, matching the
# unchanged activeTab="0" in the workbook view.
$wb.Worksheets.Item(1).Activate()

Write-Host "Sheet order:"
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $s = $wb.Worksheets.Item($i)
    Write-Host " $i : $($s.Name)"
}
